# Update the CandyStore "getOrder" webservice URL from the old EC2 host
# (54.173.199.137) to localhost, per commit "Changed 54.173.199.137 to localhost".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("WSGetOrderData")

$newUrl = "http://localhost:8080/CandyStore/candyorderREST/getOrder/"

# 1. Update the cell text itself (A2 held the old, long EC2 URL).
$ws.Range("A2").Value = $newUrl

# 2. Update the hyperlink attached to A2 so its displayed text matches the
#    new URL, while leaving the underlying hyperlink target/address as-is.
$oldAddress = "http://54.173.199.137:8080/CandyStore"
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), $oldAddress, "", "", $newUrl)

# Re-adding the hyperlink re-applies Excel's default blue/underlined
# "Hyperlink" look; restore the cell's original (non-hyperlink) font so its
# appearance/fill stay as they were before the edit.
$ws.Range("A2").Font.Underline = $false
$ws.Range("A2").Font.ColorIndex = 1
